$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update duration column (C2:C65) from 0.5 to 0.75
$ws.Range("C2:C65").Value = 0.75

# Append 8 new rest rows (66-73): name, pitch(duration constant), duration
$restNames = @("rest_1","rest_2","rest_3","rest_4","rest_5","rest_6","rest_7","rest_8")
$row = 66
foreach ($name in $restNames) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 0.0067
    $ws.Cells.Item($row, 3).Value = 0.75
    $row = $row + 1
}

# Update selection to match the author's final cursor position
$ws.Range("E65").Select()

# Reflect the author's window geometry on the active window
$win = $excel.ActiveWindow
$win.Left = 2820
$win.Top = 1220
$win.Width = 26280
$win.Height = 16820
